# Issue 769: Migrate from JUnit to TestNG (and fix QUnit casing)
#
# The "JUnit, Qunit, Selenium" caption on the testing-tools shape becomes
# "TestNG, QUnit, Selenium". Each word lives in its own <a:r> run, so we
# edit the text in place via TextRange.Characters(start, length) instead
# of rewriting the whole TextRange.Text, to leave every run's formatting
# (and the rest of the slide) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-Substring($TextRange, $OldWord, $NewWord) {
    $current = $TextRange.Text
    $idx = $current.IndexOf($OldWord)
    if ($idx -ge 0) {
        $sub = $TextRange.Characters($idx + 1, $OldWord.Length)
        $sub.Text = $NewWord
    }
}

foreach ($shp in $s.Shapes) {
    if (-not $shp.HasTextFrame) { continue }
    $tf = $shp.TextFrame
    if (-not $tf.HasText) { continue }
    $tr = $tf.TextRange

    if ($tr.Text -eq "JUnit, Qunit, Selenium") {
        Set-Substring $tr "Qunit" "QUnit"
        Set-Substring $tr "JUnit" "TestNG"
    }
}
